# Updates to NA Param file column and row names
#
# The "Capex y-intercept (" / "Capex flow coefficient (" row labels on the
# Large / Medium / Small sheets used the literal euro-sign glyph "€" inside
# their (bold) unit suffix run. Replace the glyph with the word "euros" so
# the label now reads "...(euros/t/yr/100km)" / "...(euros/t^2/yr^2/100km)"
# on every sheet, preserving the existing rich-text split: a plain prefix
# run ("Capex y-intercept (euros" / "Capex flow coefficient (euros") and a
# bold suffix run carrying the unit text.

$wb = $excel.ActiveWorkbook

$yInterceptFull   = "Capex y-intercept (euros/t/yr/100km)"
$yInterceptPrefix = "Capex y-intercept (euros"
$flowCoeffFull    = "Capex flow coefficient (euros/t^2/yr^2/100km)"
$flowCoeffPrefix  = "Capex flow coefficient (euros"

$wsLarge  = $wb.Worksheets.Item("Large")
$wsMedium = $wb.Worksheets.Item("Medium")
$wsSmall  = $wb.Worksheets.Item("Small")

# Re-text the Large sheet's two Capex labels directly, restoring the bold
# formatting on the unit suffix via Characters (this is the rich-text run
# that used to carry the euro sign).
$rngA2 = $wsLarge.Range("A2")
$rngA2.Value = $yInterceptFull
$sufA2 = $rngA2.Characters($yInterceptPrefix.Length + 1, $yInterceptFull.Length - $yInterceptPrefix.Length)
$sufA2.Font.Bold = $true

$rngA3 = $wsLarge.Range("A3")
$rngA3.Value = $flowCoeffFull
$sufA3 = $rngA3.Characters($flowCoeffPrefix.Length + 1, $flowCoeffFull.Length - $flowCoeffPrefix.Length)
$sufA3.Font.Bold = $true

# Propagate the identical re-texted + re-formatted labels to Medium and
# Small via copy/paste so both sheets end up with the same corrected text
# and rich-text formatting as Large.
$wsLarge.Range("A2:A3").Copy()
$wsMedium.Range("A2").PasteSpecial(-4104)

$wsLarge.Range("A2:A3").Copy()
$wsSmall.Range("A2").PasteSpecial(-4104)
